$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The edit described by the diff "accepts"/discards the two reviewer
# comments in the document: it removes the <w:commentRangeStart>,
# <w:commentRangeEnd> and <w:commentReference> markers (and the now
# completely empty trailing run that used to host the comment
# reference) while keeping the commented text itself, merged back
# into the flow of the paragraph it was anchored in.
# ------------------------------------------------------------------

# 1) Delete every comment in the document (walk backwards so the
#    collection indices stay valid as items are removed). This strips
#    out commentRangeStart / commentRangeEnd / commentReference and
#    the comments.xml payload.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# 2) The comment boundaries used to split what is logically one run
#    of text into two adjacent runs. Re-merge each pair back into a
#    single run by searching for the concatenation of the two pieces
#    of text and "replacing" it with itself - Word collapses the
#    match back into one run when it performs the substitution.

$loc1 = "Το σύστημα της Ολυμπίας οδού καταγράφει όλες τις διελεύσεις που γίνονται σε κάθε σταθμό τις αποθηκεύει τοπικά και τις αποστέλλει στο σύστημα, προκειμένου να γίνει ο υπολογισμός των οφειλών. Αν η αποστολή δεν επιτύχει γίνεται νέα απόπειρα μετά από καθορισμένο χρονικό διάστημα. Παράλληλα μπορεί να επισκεφθεί τον ιστότοπο και αφού συνδεθεί να λάβει συγκεκριμένες οφειλές από και προς αυτόν και να μεταβεί στο σύστημα πληρωμών ώστε να τις αποπληρώσει. "

$d.Content.Find.Execute($loc1, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $loc1, 2) | Out-Null

$loc2 = "γρήγορη απόκριση του συστήματος στην προσθήκη και στην ενημέρωση των δεδομένων, όπως για παράδειγμα η αγορά ενός νέου "

$d.Content.Find.Execute($loc2, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $loc2, 2) | Out-Null

Write-Output "Comments remaining: $($d.Comments.Count)"
